$d = $word.ActiveDocument

# Remove the leftover "_GoBack" bookmark (an empty bookmark Word leaves behind
# to mark the last edit position). Deleting it collapses the paragraph that
# only contained the bookmarkStart/bookmarkEnd pair down to an empty paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
